{"js": "// Replace the 25 \"dividend\u00f7divisor=quotient, remainder\" answer strings\n// inside the worksheet table with their regenerated values, matching the\n// author's re-run of the answer generator (commit c8c62b6). Every source\n// string occurs exactly once in the document, so an exact full-paragraph\n// text match is sufficient and avoids any partial / overlapping matches.\nconst replacements = {\n  \"682\u00f74=170, 2\": \"822\u00f74=205, 2\",\n  \"660\u00f78=82, 4\": \"626\u00f78=78, 2\",\n  \"325\u00f77=46, 3\": \"795\u00f75=159, 0\",\n  \"700\u00f72=350, 0\": \"284\u00f73=94, 2\",\n  \"260\u00f73=86, 2\": \"997\u00f77=142, 3\",\n  \"995\u00f75=199, 0\": \"211\u00f73=70, 1\",\n  \"433\u00f76=72, 1\": \"971\u00f75=194, 1\",\n  \"321\u00f75=64, 1\": \"854\u00f78=106, 6\",\n  \"511\u00f75=102, 1\": \"910\u00f75=182, 0\",\n  \"872\u00f74=218, 0\": \"270\u00f78=33, 6\",\n  \"973\u00f77=139, 0\": \"207\u00f75=41, 2\",\n  \"844\u00f78=105, 4\": \"762\u00f73=254, 0\",\n  \"362\u00f73=120, 2\": \"166\u00f72=83, 0\",\n  \"817\u00f72=408, 1\": \"894\u00f74=223, 2\",\n  \"815\u00f78=101, 7\": \"696\u00f72=348, 0\",\n  \"396\u00f72=198, 0\": \"247\u00f75=49, 2\",\n  \"117\u00f78=14, 5\": \"110\u00f76=18, 2\",\n  \"394\u00f72=197, 0\": \"281\u00f74=70, 1\",\n  \"838\u00f76=139, 4\": \"595\u00f74=148, 3\",\n  \"271\u00f79=30, 1\": \"129\u00f74=32, 1\",\n  \"113\u00f75=22, 3\": \"301\u00f73=100, 1\",\n  \"895\u00f76=149, 1\": \"716\u00f72=358, 0\",\n  \"428\u00f72=214, 0\": \"426\u00f79=47, 3\",\n  \"629\u00f78=78, 5\": \"182\u00f74=45, 2\",\n  \"228\u00f74=57, 0\": \"296\u00f78=37, 0\",\n};\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  paragraph.load(\"text\");\n}\nawait context.sync();\n\nlet replacedCount = 0;\nfor (const paragraph of paragraphs.items) {\n  const currentText = paragraph.text;\n  if (Object.prototype.hasOwnProperty.call(replacements, currentText)) {\n    paragraph.insertText(replacements[currentText], \"Replace\");\n    replacedCount++;\n  }\n}\nawait context.sync();\n\nif (replacedCount !== Object.keys(replacements).length) {\n  throw new Error(\n    `Expected to replace ${Object.keys(replacements).length} paragraphs but replaced ${replacedCount}`\n  );\n}\n", "ps1": "# Replace the 25 \"dividend\u00f7divisor=quotient, remainder\" answer strings\n# inside the worksheet table with their regenerated values, matching the\n# author's re-run of the answer generator (commit c8c62b6). Every source\n# string occurs exactly once in the document, so Find/Replace (one hit\n# each, wdReplaceAll=2 just to be safe) is sufficient and leaves run\n# formatting (font/size) untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"682\u00f74=170, 2\"; Replace = \"822\u00f74=205, 2\" },\n    @{ Find = \"660\u00f78=82, 4\"; Replace = \"626\u00f78=78, 2\" },\n    @{ Find = \"325\u00f77=46, 3\"; Replace = \"795\u00f75=159, 0\" },\n    @{ Find = \"700\u00f72=350, 0\"; Replace = \"284\u00f73=94, 2\" },\n    @{ Find = \"260\u00f73=86, 2\"; Replace = \"997\u00f77=142, 3\" },\n    @{ Find = \"995\u00f75=199, 0\"; Replace = \"211\u00f73=70, 1\" },\n    @{ Find = \"433\u00f76=72, 1\"; Replace = \"971\u00f75=194, 1\" },\n    @{ Find = \"321\u00f75=64, 1\"; Replace = \"854\u00f78=106, 6\" },\n    @{ Find = \"511\u00f75=102, 1\"; Replace = \"910\u00f75=182, 0\" },\n    @{ Find = \"872\u00f74=218, 0\"; Replace = \"270\u00f78=33, 6\" },\n    @{ Find = \"973\u00f77=139, 0\"; Replace = \"207\u00f75=41, 2\" },\n    @{ Find = \"844\u00f78=105, 4\"; Replace = \"762\u00f73=254, 0\" },\n    @{ Find = \"362\u00f73=120, 2\"; Replace = \"166\u00f72=83, 0\" },\n    @{ Find = \"817\u00f72=408, 1\"; Replace = \"894\u00f74=223, 2\" },\n    @{ Find = \"815\u00f78=101, 7\"; Replace = \"696\u00f72=348, 0\" },\n    @{ Find = \"396\u00f72=198, 0\"; Replace = \"247\u00f75=49, 2\" },\n    @{ Find = \"117\u00f78=14, 5\"; Replace = \"110\u00f76=18, 2\" },\n    @{ Find = \"394\u00f72=197, 0\"; Replace = \"281\u00f74=70, 1\" },\n    @{ Find = \"838\u00f76=139, 4\"; Replace = \"595\u00f74=148, 3\" },\n    @{ Find = \"271\u00f79=30, 1\"; Replace = \"129\u00f74=32, 1\" },\n    @{ Find = \"113\u00f75=22, 3\"; Replace = \"301\u00f73=100, 1\" },\n    @{ Find = \"895\u00f76=149, 1\"; Replace = \"716\u00f72=358, 0\" },\n    @{ Find = \"428\u00f72=214, 0\"; Replace = \"426\u00f79=47, 3\" },\n    @{ Find = \"629\u00f78=78, 5\"; Replace = \"182\u00f74=45, 2\" },\n    @{ Find = \"228\u00f74=57, 0\"; Replace = \"296\u00f78=37, 0\" }\n)\n\n$successCount = 0\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Find\n    $find.Replacement.Text = $r.Replace\n    $found = $find.Execute($r.Find, $false, $false, $false, $false, $false, $true, 1, $false, $r.Replace, 2)\n    if ($found) {\n        $successCount = $successCount + 1\n    } else {\n        Write-Output \"WARNING: text not found for replacement: $($r.Find)\"\n    }\n}\n\nif ($successCount -ne $replacements.Count) {\n    throw \"Expected to replace $($replacements.Count) strings but only replaced $successCount\"\n}\n\nWrite-Output \"done: replaced $successCount strings\"\n"}
